{"js": "// The document contains a plain-text content control (w:sdt) whose\n// alias/tag is \"M\u00fcd\u00fcr\", wrapping a single centered paragraph that holds\n// one space character. The edit removes the content control wrapper\n// while keeping its paragraph content in place (i.e. \"Remove Content\n// Control\" in the Word UI / ContentControl.Delete(keepContent:=True)).\n//\n// NOTE: `document.contentControls` (the property) only surfaces\n// rich-text SDTs; this one is a plain-text SDT (it has a <w:text/>\n// marker), so it must be located via `getContentControls()` (which\n// enumerates every SDT type) instead.\nconst contentControls = context.document.getContentControls();\nconst target = contentControls.getByTag(\"M\u00fcd\u00fcr\");\ntarget.load(\"items\");\nawait context.sync();\n\nif (target.items.length > 0) {\n  // keepContent = true: unwrap the control, leaving its paragraph/run behind.\n  target.items[0].delete(true);\n  await context.sync();\n}\n", "ps1": "# The document contains a plain-text content control (Word \"legacy\"\n# Building Block SDT) whose Title/Tag is \"M\u00fcd\u00fcr\", wrapping a single\n# centered paragraph that holds one space character. The edit removes\n# the content control wrapper while keeping its paragraph content in\n# place - i.e. the same thing as selecting the control in the Word UI\n# and choosing \"Remove Content Control\" (ContentControl.Delete with\n# DeleteContents = False).\n$d = $word.ActiveDocument\n\nforeach ($cc in $d.ContentControls) {\n    if ($cc.Tag -eq \"M\u00fcd\u00fcr\") {\n        # False => keep the control's contents, only remove the wrapper.\n        $cc.Delete($false)\n    }\n}\n"}
